$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, taken from the refreshed cryptos listing.
# Values must land in the sheet as plain TEXT (the source column is
# text-formatted, e.g. "56.206.36" / "3.229.08" with thousand-grouping
# dots, and the percent column keeps its padding spaces), so every write
# below is guarded to prevent Excel/COM from reinterpreting a
# numeric-looking string (e.g. "397.30", "3.32") as a Double and silently
# dropping the trailing zero / formatting.
$updates = @(
    @{ Cell = "D2"; Value = "56.206.36" },
    @{ Cell = "E2"; Value = "  +9.24%  " },
    @{ Cell = "D3"; Value = "3.229.08" },
    @{ Cell = "E3"; Value = "  +3.95%  " },
    @{ Cell = "E4"; Value = "  +0.02%  " },
    @{ Cell = "D5"; Value = "397.30" },
    @{ Cell = "E5"; Value = "  +2.13%  " },
    @{ Cell = "D6"; Value = "110.86" },
    @{ Cell = "E6"; Value = "  +6.56%  " },
    @{ Cell = "E7"; Value = "  +2.59%  " },
    @{ Cell = "E9"; Value = "  +5.11%  " },
    @{ Cell = "D10"; Value = "39.31" },
    @{ Cell = "E10"; Value = "  +5.77%  " },
    @{ Cell = "D11"; Value = "0.0913" },
    @{ Cell = "E11"; Value = "  +6.31%  " },
    @{ Cell = "E12"; Value = "  +2.14%  " },
    @{ Cell = "D13"; Value = "3.736.54" },
    @{ Cell = "E13"; Value = "  +3.81%  " },
    @{ Cell = "D14"; Value = "8.07" },
    @{ Cell = "E14"; Value = "  +3.53%  " },
    @{ Cell = "D15"; Value = "19.05" },
    @{ Cell = "E15"; Value = "  +2.54%  " },
    @{ Cell = "D16"; Value = "3.229.86" },
    @{ Cell = "E16"; Value = "  +3.90%  " },
    @{ Cell = "E17"; Value = "  +4.92%  " },
    @{ Cell = "D18"; Value = "10.75" },
    @{ Cell = "E18"; Value = "  +1.13%  " },
    @{ Cell = "D19"; Value = "56.046.35" },
    @{ Cell = "E19"; Value = "  +8.70%  " },
    @{ Cell = "D20"; Value = "3.32" },
    @{ Cell = "E20"; Value = "  +1.53%  " },
    @{ Cell = "E21"; Value = "  +5.92%  " },
    @{ Cell = "D22"; Value = "13.01" },
    @{ Cell = "E22"; Value = "  +3.93%  " },
    @{ Cell = "D23"; Value = "299.21" },
    @{ Cell = "E23"; Value = "  +12.12%  " },
    @{ Cell = "D24"; Value = "75.41" },
    @{ Cell = "E24"; Value = "  +7.28%  " },
    @{ Cell = "E25"; Value = "  +1.40%  " },
    @{ Cell = "D26"; Value = "8.19" },
    @{ Cell = "E26"; Value = "  +1.41%  " },
    @{ Cell = "D27"; Value = "28.20" },
    @{ Cell = "E27"; Value = "  +2.93%  " },
    @{ Cell = "D28"; Value = "7.50" },
    @{ Cell = "E28"; Value = "  +4.26%  " },
    @{ Cell = "D29"; Value = "0.172" },
    @{ Cell = "E29"; Value = "  +4.17%  " },
    @{ Cell = "E30"; Value = "  +0.44%  " },
    @{ Cell = "E31"; Value = "  +3.43%  " },
    @{ Cell = "E32"; Value = "  +6.70%  " },
    @{ Cell = "D33"; Value = "0.0492" },
    @{ Cell = "E33"; Value = "  +3.33%  " },
    @{ Cell = "D34"; Value = "36.22" },
    @{ Cell = "E34"; Value = "  +1.03%  " },
    @{ Cell = "E35"; Value = "  +2.67%  " },
    @{ Cell = "D36"; Value = "51.29" },
    @{ Cell = "E36"; Value = "  +2.36%  " },
    @{ Cell = "E37"; Value = "  +25.53%  " },
    @{ Cell = "E38"; Value = "  +3.87%  " },
    @{ Cell = "D39"; Value = "0.999" },
    @{ Cell = "E39"; Value = "  +0.02%  " },
    @{ Cell = "B40"; Value = "Monero" },
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" },
    @{ Cell = "D40"; Value = "135.05" },
    @{ Cell = "E40"; Value = "  +4.31%  " },
    @{ Cell = "B41"; Value = "ARBITRUM" },
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb" },
    @{ Cell = "D41"; Value = "1.92" },
    @{ Cell = "E41"; Value = "  +3.00%  " },
    @{ Cell = "B42"; Value = "Celestia" },
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia" },
    @{ Cell = "D42"; Value = "17.30" },
    @{ Cell = "E42"; Value = "  +4.19%  " },
    @{ Cell = "B43"; Value = "NEARProtocol" },
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near" },
    @{ Cell = "D43"; Value = "4.01" },
    @{ Cell = "E43"; Value = "  +4.92%  " },
    @{ Cell = "E44"; Value = "  +2.86%  " },
    @{ Cell = "D45"; Value = "0.283" },
    @{ Cell = "E45"; Value = "  -2.81%  " },
    @{ Cell = "D46"; Value = "22.28" },
    @{ Cell = "E46"; Value = "  +0.48%  " },
    @{ Cell = "D47"; Value = "2.13" },
    @{ Cell = "E47"; Value = "  +47.00%  " },
    @{ Cell = "E48"; Value = "  +1.77%  " },
    @{ Cell = "E49"; Value = "  -0.83%  " },
    @{ Cell = "D50"; Value = "2.133.08" },
    @{ Cell = "E50"; Value = "  +2.54%  " },
    @{ Cell = "D51"; Value = "0.0363" },
    @{ Cell = "E51"; Value = "  +9.38%  " }
)

foreach ($update in $updates) {
    $cellRef = $update.Cell
    $newVal = $update.Value
    $rng = $ws.Range($cellRef)

    # Plain decimal-looking text (optional sign, digits, optional single
    # decimal point) is what Excel's COM layer will happily reinterpret as
    # a Double on assignment. Thousand-dotted prices ("56.206.36"), percent
    # strings, and coin names/URLs never match this and go through as-is.
    $looksNumeric = ($newVal -match '^[+-]?(\d+\.?\d*|\.\d+)$')

    if ($looksNumeric) {
        # Force text storage so "397.30" doesn't become the number 397.3.
        $rng.NumberFormat = "@"
        $rng.Value = $newVal
        # Drop back to the sheet's default style so no stray number-format
        # style is left attached to the cell.
        $rng.Style = "Normal"
    } else {
        $rng.Value = $newVal
    }
}
